$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Step 1: insert 7 new rows before (old) row 11, to make room for the 7 new rows of problems
$anchorRow = $t.Rows.Item(11)
for ($i = 0; $i -lt 7; $i++) {
    $newRow = $t.Rows.Add($anchorRow)
}

# Step 2: remove the trailing 7 rows (the table now has 27 rows; drop the last 7
# so we end up with 20 rows total, matching the final layout)
for ($i = 0; $i -lt 7; $i++) {
    $t.Rows.Item($t.Rows.Count).Delete()
}

Write-Host "Row count:" $t.Rows.Count

# Step 3: write every cell value for the final 20x5 grid of problems
$t.Cell(1,1).Range.Text = "91-74=17"
$t.Cell(1,2).Range.Text = "66+27=93"
$t.Cell(1,3).Range.Text = "5+28=33"
$t.Cell(1,4).Range.Text = "47+39=86"
$t.Cell(1,5).Range.Text = "19+39=58"
$t.Cell(2,1).Range.Text = "16+8=24"
$t.Cell(2,2).Range.Text = "44+27=71"
$t.Cell(2,3).Range.Text = "52-37=15"
$t.Cell(2,4).Range.Text = "52-29=23"
$t.Cell(2,5).Range.Text = "46-7=39"
$t.Cell(3,1).Range.Text = "4+18=22"
$t.Cell(3,2).Range.Text = "19+79=98"
$t.Cell(3,3).Range.Text = "58+3=61"
$t.Cell(3,4).Range.Text = "34-15=19"
$t.Cell(3,5).Range.Text = "91-9=82"
$t.Cell(4,1).Range.Text = "42-38=4"
$t.Cell(4,2).Range.Text = "45-36=9"
$t.Cell(4,3).Range.Text = "32+49=81"
$t.Cell(4,4).Range.Text = "16+39=55"
$t.Cell(4,5).Range.Text = "92-88=4"
$t.Cell(5,1).Range.Text = "92-68=24"
$t.Cell(5,2).Range.Text = "14+17=31"
$t.Cell(5,3).Range.Text = "28+47=75"
$t.Cell(5,4).Range.Text = "72-48=24"
$t.Cell(5,5).Range.Text = "74-8=66"
$t.Cell(6,1).Range.Text = "3+28=31"
$t.Cell(6,2).Range.Text = "57+35=92"
$t.Cell(6,3).Range.Text = "18+7=25"
$t.Cell(6,4).Range.Text = "58-49=9"
$t.Cell(6,5).Range.Text = "60-2=58"
$t.Cell(7,1).Range.Text = "27+34=61"
$t.Cell(7,2).Range.Text = "36+45=81"
$t.Cell(7,3).Range.Text = "73-8=65"
$t.Cell(7,4).Range.Text = "88+3=91"
$t.Cell(7,5).Range.Text = "53-39=14"
$t.Cell(8,1).Range.Text = "49+23=72"
$t.Cell(8,2).Range.Text = "28-9=19"
$t.Cell(8,3).Range.Text = "64-49=15"
$t.Cell(8,4).Range.Text = "60-23=37"
$t.Cell(8,5).Range.Text = "75-36=39"
$t.Cell(9,1).Range.Text = "19+43=62"
$t.Cell(9,2).Range.Text = "18+59=77"
$t.Cell(9,3).Range.Text = "82-69=13"
$t.Cell(9,4).Range.Text = "83-35=48"
$t.Cell(9,5).Range.Text = "32-4=28"
$t.Cell(10,1).Range.Text = "35+39=74"
$t.Cell(10,2).Range.Text = "6+5=11"
$t.Cell(10,3).Range.Text = "57-29=28"
$t.Cell(10,4).Range.Text = "18+77=95"
$t.Cell(10,5).Range.Text = "56+37=93"
$t.Cell(11,1).Range.Text = "71-4=67"
$t.Cell(11,2).Range.Text = "19+45=64"
$t.Cell(11,3).Range.Text = "26+9=35"
$t.Cell(11,4).Range.Text = "82-79=3"
$t.Cell(11,5).Range.Text = "90-55=35"
$t.Cell(12,1).Range.Text = "84+8=92"
$t.Cell(12,2).Range.Text = "93-44=49"
$t.Cell(12,3).Range.Text = "8+68=76"
$t.Cell(12,4).Range.Text = "78+16=94"
$t.Cell(12,5).Range.Text = "51-42=9"
$t.Cell(13,1).Range.Text = "6+8=14"
$t.Cell(13,2).Range.Text = "63-56=7"
$t.Cell(13,3).Range.Text = "36-28=8"
$t.Cell(13,4).Range.Text = "13+48=61"
$t.Cell(13,5).Range.Text = "29+53=82"
$t.Cell(14,1).Range.Text = "47+17=64"
$t.Cell(14,2).Range.Text = "34-16=18"
$t.Cell(14,3).Range.Text = "67-39=28"
$t.Cell(14,4).Range.Text = "88+5=93"
$t.Cell(14,5).Range.Text = "42-36=6"
$t.Cell(15,1).Range.Text = "16+78=94"
$t.Cell(15,2).Range.Text = "60-29=31"
$t.Cell(15,3).Range.Text = "48+18=66"
$t.Cell(15,4).Range.Text = "87-39=48"
$t.Cell(15,5).Range.Text = "6+7=13"
$t.Cell(16,1).Range.Text = "90-75=15"
$t.Cell(16,2).Range.Text = "18+43=61"
$t.Cell(16,3).Range.Text = "48+28=76"
$t.Cell(16,4).Range.Text = "65-47=18"
$t.Cell(16,5).Range.Text = "56-19=37"
$t.Cell(17,1).Range.Text = "55+8=63"
$t.Cell(17,2).Range.Text = "76-19=57"
$t.Cell(17,3).Range.Text = "7+18=25"
$t.Cell(17,4).Range.Text = "60-24=36"
$t.Cell(17,5).Range.Text = "46+18=64"
$t.Cell(18,1).Range.Text = "69+26=95"
$t.Cell(18,2).Range.Text = "4+88=92"
$t.Cell(18,3).Range.Text = "38+25=63"
$t.Cell(18,4).Range.Text = "20-9=11"
$t.Cell(18,5).Range.Text = "47+16=63"
$t.Cell(19,1).Range.Text = "60-41=19"
$t.Cell(19,2).Range.Text = "82-49=33"
$t.Cell(19,3).Range.Text = "83-7=76"
$t.Cell(19,4).Range.Text = "90-19=71"
$t.Cell(19,5).Range.Text = "23-14=9"
$t.Cell(20,1).Range.Text = "57-29=28"
$t.Cell(20,2).Range.Text = "55-39=16"
$t.Cell(20,3).Range.Text = "20-16=4"
$t.Cell(20,4).Range.Text = "28+54=82"
$t.Cell(20,5).Range.Text = "35+47=82"
